# Update "Pagos" (F) and "Inscrições homologadas" (H) counts for a few rows.
# For each affected row, F increases by 1 and H (= F + G) increases by 1 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F36").Value = 88
$ws.Range("H36").Value = 120

$ws.Range("F41").Value = 35
$ws.Range("H41").Value = 46

$ws.Range("F72").Value = 42
$ws.Range("H72").Value = 53

$ws.Range("F74").Value = 12
$ws.Range("H74").Value = 16
